$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204; this shifts existing rows 204..291 down to 205..292
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new record's data
$ws.Cells.Item(204, 1).Value = 3
$ws.Cells.Item(204, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(204, 3).Value = "Coquimbo"
$ws.Cells.Item(204, 4).Value = 44704
$ws.Cells.Item(204, 5).Value = 5
$ws.Cells.Item(204, 6).Value = 100112001
$ws.Cells.Item(204, 7).Value = "Berenjena"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 105
$ws.Cells.Item(204, 11).Value = 6000
$ws.Cells.Item(204, 12).Value = 7000
$ws.Cells.Item(204, 13).Value = 6476
$ws.Cells.Item(204, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(204, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(204, 16).Value = 108
$ws.Cells.Item(204, 17).Value = 60
$ws.Cells.Item(204, 18).Value = "Hortaliza"
